# Generate Report for Handback
#
# The handback transform failed for the "3fe9a775-3405-45dc-b710-a28b5256dd89.md"
# file in both the zh-cn and de-de locales. Update the status on the Overview
# sheet and record the error detail on each locale sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Row 3 everywhere is the "3fe9a775-3405-45dc-b710-a28b5256dd89.md" file.
# Its Status - mirrored on the Overview sheet (zh-cn column E, de-de column
# F) and on each locale sheet (column C) - flips from "Ready for handoff" to
# "Handback transform failed".
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# zh-cn / de-de sheets: row 3 is the same file. Populate the "Error Detail"
# column (P) with the handback/handoff file-name mismatch message.
$wsZhCn.Range("P3").Value = "Handback file name: eieo2tvr.bq0 is different with handoff file name: 3fe9a775-3405-45dc-b710-a28b5256dd89.fd08e891fbb9f9fa7e7ce17f54a6a997d232bb8c.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: eieo2tvr.bq0 is different with handoff file name: 3fe9a775-3405-45dc-b710-a28b5256dd89.fd08e891fbb9f9fa7e7ce17f54a6a997d232bb8c.de-de."

# Widen column P (Error Detail) on both locale sheets so the message is
# readable - matches an XML column width of 40 (ColumnWidth is offset by the
# default 5/6-character padding Excel adds on top of the stored value).
$wsZhCn.Columns.Item(16).ColumnWidth = 40 - 5/6
$wsDeDe.Columns.Item(16).ColumnWidth = 40 - 5/6
